# Delete the last slide (slide 11, "세부 사항" / access-policy detail slide).
# PowerPoint will renumber the remaining relationship ids (notesMaster
# moves from rId13 -> rId12, etc.) and drop the corresponding notes
# slide (notesSlide9.xml) automatically as part of removing the slide.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)
$s.Delete()
